# Remove the "tại HỆ THỐNG" rows from the "Lương" (Salary) sheet.
# This deletes the 7 detail rows (Chiết khấu/Đơn/Công phụ phẫu/Ứng lương
# "tại HỆ THỐNG") as well as the "Tổng lương tại HỆ THỐNG" summary row,
# shifting all subsequent rows up so the remaining data stays contiguous.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

# Delete the block of 7 "... tại HỆ THỐNG" detail rows (originally rows 4-10).
$ws.Range("A4:A10").EntireRow.Delete()

# After the block above is removed, the "Tổng lương tại HỆ THỐNG" row
# (originally row 35) has shifted up to row 28. Delete it too.
$ws.Range("A28").EntireRow.Delete()
